$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column values that look numeric (e.g. "1.00", "58.80") keep their
# exact textual formatting instead of being auto-converted to numbers by Excel.
$dCells = @("D2","D3","D4","D5","D6","D8","D11","D14","D15","D16","D17","D19","D20","D21","D22","D23","D25","D26","D29","D30","D31","D32","D34","D36","D37","D39","D41","D42","D43","D45","D46","D47","D48","D49","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "54.342.46"
$ws.Range("E2").Value = "  -7.22%  "
$ws.Range("D3").Value = "2.871.02"
$ws.Range("E3").Value = "  -10.03%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "470.92"
$ws.Range("E5").Value = "  -11.48%  "
$ws.Range("D6").Value = "125.72"
$ws.Range("E6").Value = "  -6.61%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "2.867.87"
$ws.Range("E8").Value = "  -10.17%  "
$ws.Range("E9").Value = "  -11.78%  "
$ws.Range("E10").Value = "  -10.83%  "
$ws.Range("D11").Value = "0.0959"
$ws.Range("E11").Value = "  -14.98%  "
$ws.Range("E12").Value = "  -15.71%  "
$ws.Range("E13").Value = "  -4.54%  "
$ws.Range("D14").Value = "3.360.45"
$ws.Range("E14").Value = "  -10.19%  "
$ws.Range("D15").Value = "23.25"
$ws.Range("E15").Value = "  -9.52%  "
$ws.Range("D16").Value = "54.337.03"
$ws.Range("E16").Value = "  -7.37%  "
$ws.Range("D17").Value = "2.867.56"
$ws.Range("E17").Value = "  -10.33%  "
$ws.Range("E18").Value = "  -14.13%  "
$ws.Range("D19").Value = "5.33"
$ws.Range("E19").Value = "  -9.46%  "
$ws.Range("D20").Value = "11.43"
$ws.Range("E20").Value = "  -13.20%  "
$ws.Range("D21").Value = "7.04"
$ws.Range("E21").Value = "  -13.04%  "
$ws.Range("D22").Value = "294.95"
$ws.Range("E22").Value = "  -17.60%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  -13.77%  "
$ws.Range("D25").Value = "58.80"
$ws.Range("E25").Value = "  -15.58%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -9.76%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "0.0₃0805"
$ws.Range("E29").Value = "  -15.21%  "
$ws.Range("D30").Value = "6.17"
$ws.Range("E30").Value = "  -12.19%  "
$ws.Range("D31").Value = "1.13"
$ws.Range("E31").Value = "  -5.24%  "
$ws.Range("D32").Value = "6.20"
$ws.Range("E32").Value = "  -11.30%  "
$ws.Range("E33").Value = "  -15.62%  "
$ws.Range("D34").Value = "18.84"
$ws.Range("E34").Value = "  -12.66%  "
$ws.Range("E35").Value = "  -13.47%  "
$ws.Range("D36").Value = "135.17"
$ws.Range("E36").Value = "  -16.22%  "
$ws.Range("D37").Value = "5.42"
$ws.Range("E37").Value = "  -14.11%  "
$ws.Range("E38").Value = "  -14.05%  "
$ws.Range("D39").Value = "23.11"
$ws.Range("E39").Value = "  -10.19%  "
$ws.Range("E40").Value = "  -12.03%  "
$ws.Range("D41").Value = "2.894.95"
$ws.Range("E41").Value = "  -10.14%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "34.94"
$ws.Range("E43").Value = "  -14.33%  "
$ws.Range("E44").Value = "  -12.94%  "
$ws.Range("D45").Value = "0.601"
$ws.Range("E45").Value = "  -15.31%  "
$ws.Range("D46").Value = "1.31"
$ws.Range("E46").Value = "  -11.47%  "
$ws.Range("D47").Value = "3.38"
$ws.Range("E47").Value = "  -15.22%  "
$ws.Range("D48").Value = "2.049.61"
$ws.Range("E48").Value = "  -10.27%  "
$ws.Range("D49").Value = "5.36"
$ws.Range("E49").Value = "  -13.94%  "
$ws.Range("D50").Value = "17.92"
$ws.Range("E50").Value = "  -12.19%  "
$ws.Range("E51").Value = "  -11.12%  "
